# Task: rename physical-register labels "pNN" -> "rNN" throughout the
# pipeline-renaming table, and clear out the G11:I11 allocation (the row
# whose branch outcome is MISS), shifting the "p39" that used to sit at
# G12 down to become "r38" since one fewer physical register is now live.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G column (destination physical register) ---
$ws.Range("G4").Value  = "r32"
$ws.Range("G5").Value  = "r33"
$ws.Range("G6").Value  = "r34"
$ws.Range("G7").Value  = "r35"
$ws.Range("G8").Value  = "r36"
$ws.Range("G9").Value  = "r37"
$ws.Range("G10").Value = "r38"
$ws.Range("G12").Value = "r38"
$ws.Range("G13").Value = "r39"

# --- H column (rs1, ready) ---
$ws.Range("H5").Value  = "r32,YES"
$ws.Range("H6").Value  = "r32,YES"
$ws.Range("H7").Value  = "r33,YES"
$ws.Range("H8").Value  = "r35,NO"
$ws.Range("H9").Value  = "r36,NO"
$ws.Range("H10").Value = "r36,YES"
$ws.Range("H12").Value = "r36,YES"
$ws.Range("H13").Value = "r32,YES"
$ws.Range("H15").Value = "r39,YES"

# --- I column (rs2, ready) ---
$ws.Range("I6").Value  = "r33,NO"
$ws.Range("I7").Value  = "r34,NO"
$ws.Range("I8").Value  = "r34,YES"
$ws.Range("I9").Value  = "r33,YES"
$ws.Range("I10").Value = "r37,NO"
$ws.Range("I12").Value = "r37,YES"
$ws.Range("I13").Value = "r33,YES"
$ws.Range("I15").Value = "r0,YES"

# Row 11 (the MISS row) loses its physical-register allocation entirely.
$ws.Range("G11:I11").ClearContents()

# Move the active selection the same way the original author's last click did.
$ws.Range("N15").Select()
